# DemoDaHuaImport1.xlsx - "xog tam phan upload file"
#
# Sheet1 originally listed 7 data rows (r2:r8). The edit:
#   - corrects the values in row 6 (B6/C6), which changes the computed
#     result of the shared formula in D6
#   - removes the last two data rows (7 and 8) entirely, shrinking the
#     used range from A1:D8 down to A1:D6
#   - leaves the cursor/selection on F7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Hang Nhap / Hang Ra values change (4124/124 -> 123/4).
# D6 = B6 - C6 is a formula and recalculates automatically to 119.
$ws.Range("B6").Value = 123
$ws.Range("C6").Value = 4

# Drop rows 7 and 8 completely (not just clear them) so the sheet's
# used range/dimension becomes A1:D6.
$ws.Rows("7:8").Delete()

# Match the saved selection/active cell.
$ws.Range("F7").Select()
